$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = 5.241599999999997
$ws.Range("A9").Value = -21.9065
$ws.Range("B12").Value = 5.358599999999996
$ws.Range("A13").Value = -22.13479999999999
$ws.Range("B14").Value = 6.064900000000002
$ws.Range("A16").Value = -21.40229999999997
$ws.Range("A18").Value = -22.18960000000001
$ws.Range("B19").Value = 9.060900000000007
$ws.Range("A20").Value = -19.61479999999998
$ws.Range("A26").Value = -21.03789999999997
$ws.Range("B26").Value = 4.262200000000004
$ws.Range("A27").Value = -21.53119999999998
$ws.Range("B27").Value = 4.904500000000004
$ws.Range("A29").Value = -21.74490000000001
$ws.Range("B29").Value = 5.415400000000001
$ws.Range("A35").Value = -19.41429999999999
$ws.Range("A36").Value = -19.81349999999999
$ws.Range("B37").Value = 8.455400000000003
$ws.Range("B38").Value = 4.356300000000001
$ws.Range("A45").Value = -21.47349999999999
$ws.Range("B47").Value = 5.315599999999998
$ws.Range("B51").Value = 5.743699999999999
$ws.Range("B52").Value = 5.639899999999999
$ws.Range("A55").Value = -22.28879999999999
$ws.Range("B55").Value = 4.873499999999996
$ws.Range("A57").Value = -22.13280000000001
$ws.Range("A69").Value = -21.7093
$ws.Range("B69").Value = 5.554199999999996
$ws.Range("B70").Value = 6.2883
$ws.Range("A76").Value = -21.9522
$ws.Range("B76").Value = 5.459899999999999
$ws.Range("A78").Value = -19.86109999999999
$ws.Range("B81").Value = 5.098700000000001
$ws.Range("A82").Value = -22.04530000000001
$ws.Range("A83").Value = -21.94009999999999
$ws.Range("B83").Value = 6.0137
$ws.Range("A93").Value = -20.66599999999998
$ws.Range("B94").Value = 5.555599999999999
$ws.Range("A97").Value = -21.91120000000002
$ws.Range("B100").Value = 5.276499999999998
$ws.Range("B102").Value = 8.516700000000005
